{"js": "// Insert a new bulleted list item \"Classifica\u00e7\u00e3o em tiles de baixa, m\u00e9dia e\n// alta complexidade\" right after the existing \"Classifica\u00e7\u00e3o com e sem\n// hiperpar\u00e2metros\" list item (same ListParagraph style / numId=3 bullet\n// list), matching the target diff.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the anchor paragraph by its exact text.\nconst anchorText = \"Classifica\u00e7\u00e3o com e sem hiperpar\u00e2metros\";\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === anchorText) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error(\"Could not locate anchor paragraph: \" + anchorText);\n}\n\n// insertParagraph on the anchor clones the anchor's paragraph formatting\n// (style, numbering/list membership, alignment, run formatting) for the new\n// paragraph, which is exactly what's needed here.\nconst newPara = anchor.insertParagraph(\n  \"Classifica\u00e7\u00e3o em tiles de baixa, m\u00e9dia e alta complexidade\",\n  Word.InsertLocation.after\n);\n\nawait context.sync();\n", "ps1": "# Insert a new bulleted list item \"Classifica\u00e7\u00e3o em tiles de baixa, m\u00e9dia e\n# alta complexidade\" right after the existing \"Classifica\u00e7\u00e3o com e sem\n# hiperpar\u00e2metros\" list item (same ListParagraph style / numId=3 bullet\n# list), matching the target diff.\n\n$d = $word.ActiveDocument\n\n$anchorText = \"Classifica\u00e7\u00e3o com e sem hiperpar\u00e2metros\"\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Trim() -eq $anchorText) {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not locate anchor paragraph: $anchorText\"\n}\n\n# Inserting an empty paragraph right after the anchor clones the anchor's\n# paragraph formatting (style, numbering/list membership, alignment, run\n# formatting), same as Word does when you press Enter at the end of a list\n# item.\n$target.Range.InsertParagraphAfter()\n\n$newPara = $target.Next()\n$newPara.Range.Text = \"Classifica\u00e7\u00e3o em tiles de baixa, m\u00e9dia e alta complexidade\"\n"}
